$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.855.19'
$ws.Range('E2').Value = '  -0.77%  '
$ws.Range('D3').Value = '1.840.40'
$ws.Range('E3').Value = '  +1.60%  '
$ws.Range('D5').Value = "'231.69"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.52%  '
$ws.Range('D6').Value = "'0.618"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.02%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = "'39.78"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.91%  '
$ws.Range('E9').Value = '  +1.95%  '
$ws.Range('D10').Value = "'0.0687"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.18%  '
$ws.Range('D11').Value = "'0.0981"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.50%  '
$ws.Range('E12').Value = '  +1.57%  '
$ws.Range('E13').Value = '  +3.18%  '
$ws.Range('D14').Value = '1.840.38'
$ws.Range('E14').Value = '  +1.51%  '
$ws.Range('D15').Value = "'0.673"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.30%  '
$ws.Range('D16').Value = "'4.65"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.08%  '
$ws.Range('D17').Value = '34.864.19'
$ws.Range('E17').Value = '  -0.57%  '
$ws.Range('D18').Value = "'69.86"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.07%  '
$ws.Range('D19').Value = '0.0₃0788'
$ws.Range('E19').Value = '  -0.47%  '
$ws.Range('D20').Value = "'240.59"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.90%  '
$ws.Range('D21').Value = "'12.18"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.05%  '
$ws.Range('D22').Value = "'4.69"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.19%  '
$ws.Range('E23').Value = '  +0.10%  '
$ws.Range('D24').Value = "'2.28"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.19%  '
$ws.Range('D25').Value = "'171.61"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.33%  '
$ws.Range('E26').Value = '  -0.82%  '
$ws.Range('D27').Value = "'17.46"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.40%  '
$ws.Range('E28').Value = '  +1.98%  '
$ws.Range('E29').Value = '  -5.11%  '
$ws.Range('E30').Value = '  +0.07%  '
$ws.Range('D31').Value = "'0.0551"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.58%  '
$ws.Range('E33').Value = '  -2.01%  '
$ws.Range('E34').Value = '  +6.66%  '
$ws.Range('E35').Value = '  +6.82%  '
$ws.Range('D36').Value = "'1.44"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +12.81%  '
$ws.Range('E37').Value = '  +1.94%  '
$ws.Range('E38').Value = '  +6.51%  '
$ws.Range('D39').Value = "'90.57"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.74%  '
$ws.Range('D40').Value = '1.345.73'
$ws.Range('E40').Value = '  +2.40%  '
$ws.Range('E41').Value = '  +0.01%  '
$ws.Range('D42').Value = "'14.92"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.62%  '
$ws.Range('D43').Value = "'2.30"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.52%  '
$ws.Range('E44').Value = '  -2.69%  '
$ws.Range('D45').Value = "'2.75"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.32%  '
$ws.Range('E46').Value = '  -0.67%  '
$ws.Range('E47').Value = '  +2.12%  '
$ws.Range('D48').Value = '2.018.97'
$ws.Range('E48').Value = '  +1.51%  '
$ws.Range('D49').Value = "'3.45"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +21.86%  '
$ws.Range('E50').Value = '  +0.09%  '
$ws.Range('E51').Value = '  +1.70%  '
